$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Empezaremos por las menos:" -> split into 3 runs with a
# proofErr (grammar-check) marker wrapped around "las".
# ---------------------------------------------------------------------------
$target1 = $d.Content
$found1 = $target1.Find.Execute("Empezaremos por las menos:", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find 'Empezaremos por las menos:' text"
}
$r1 = $d.Range($target1.Start, $target1.End)

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve">Empezaremos por </w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>las</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t xml:space="preserve"> menos:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Change 2: after " (¿?)" at the very end of the document, close off that
# paragraph with a final ".", then append a large amount of new content
# about the <build> tag's children and the Maven build lifecycle. The
# "_GoBack" bookmark that used to sit at the very end of the document now
# wraps the word "site" instead.
#
# The existing "_GoBack" bookmark occupies the zero-width slot right at the
# very end of the document (after the final paragraph's text, before its
# paragraph mark). We remove it first so our replacement text can cleanly
# take over that paragraph's tail, then the new content below recreates a
# "_GoBack" bookmark around the word "site" as the diff requires.
# ---------------------------------------------------------------------------
$goBack = $null
foreach ($bm in $d.Bookmarks) {
    if ($bm.Name -eq "_GoBack") { $goBack = $bm }
}
if ($goBack -ne $null) {
    $goBack.Delete()
}

$target2 = $d.Content
$found2 = $target2.Find.Execute(" (¿?)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find ' (¿?)' text"
}
$tailStart = $target2.End
$r2 = $d.Range($tailStart, $tailStart)

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t>.</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>
<w:r><w:t>La etiqueta &lt;build&gt; puede contener además, para dar cabida a todas esas características configurables, las siguientes:</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:t>&lt;</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>plugins</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:tab/><w:t>&lt;</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>executions</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/></w:pPr>
<w:r><w:tab/></w:r>
<w:r><w:tab/></w:r>
<w:r><w:tab/><w:t>&lt;</w:t></w:r>
<w:proofErr w:type="gramStart"/>
<w:r><w:t>configuration</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>&gt;</w:t></w:r>
</w:p>
<w:p/>
<w:p><w:pPr><w:rPr><w:b/><w:u w:val="single"/></w:rPr></w:pPr>
<w:r><w:rPr><w:b/><w:u w:val="single"/></w:rPr><w:t>Maven build lifecycle:</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t xml:space="preserve">El ciclo de vida de una aplicación </w:t></w:r>
<w:r><w:t xml:space="preserve">(el proceso de construcción y distribución de un artefacto) </w:t></w:r>
<w:r><w:t>Maven se compone de fases</w:t></w:r>
<w:r><w:t xml:space="preserve">. Hay tres ciclos de vida que vienen con el propio Maven (aunque se pueden definir otros), que son el ciclo “default”, “clean” y “site”. </w:t></w:r>
</w:p>
<w:p>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:b/></w:rPr><w:t>default</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>:</w:t></w:r>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Maneja el despliegue del proyecto</w:t></w:r>
<w:r><w:t>, y tiene las siguientes fases:</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>validate</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>compile</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>test</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>package</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>verify</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>install</w:t></w:r>
</w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>
<w:r><w:t>deploy</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:b/></w:rPr><w:t>clean</w:t></w:r>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>:</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Maneja la limpieza del proyecto</w:t></w:r>
<w:r><w:t>.</w:t></w:r>
</w:p>
<w:p/>
<w:p>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:proofErr w:type="gramStart"/>
<w:r><w:rPr><w:b/></w:rPr><w:t>site</w:t></w:r>
<w:bookmarkEnd w:id="0"/>
<w:proofErr w:type="gramEnd"/>
<w:r><w:t>:</w:t></w:r>
</w:p>
<w:p>
<w:r><w:t>Maneja la creación de documentación para el proyecto.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$r2.InsertXML($xml2)

Write-Output "Edit complete."
